# Atualização de bases das ligas, do dia: 18-05-2024 às 14:13
#
# The source feed re-pulled this league's fixtures and a handful of
# rows that are adjacent duplicates (same id block / same kickoff date)
# came back in swapped order. This swaps the full data payload
# (columns B..AB - everything except the running index in column A)
# between each such pair of rows, leaving the row index itself in
# column A untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$RowA,
        [int]$RowB
    )

    # Columns B:AB hold all of the match data (id, teams, odds, ...).
    # Column A is just the sequential row index and must stay put.
    $rangeA = $ws.Range("B$RowA`:AB$RowA")
    $rangeB = $ws.Range("B$RowB`:AB$RowB")

    # Snapshot both rows before writing anything back, so the write to
    # one row never clobbers the value we still need to read from the
    # other.
    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# Row pairs whose B:AB payloads were swapped in this update.
Swap-RowData 9 10
Swap-RowData 36 37
Swap-RowData 49 50
Swap-RowData 76 77
Swap-RowData 87 88
Swap-RowData 99 100
Swap-RowData 177 178
